$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.045.38'

$ws.Range("D3").Value = '3.209.10'
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("E4").Value = '  -0.05%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '603.32'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +4.44%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '154.06'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.98%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.209.01'
$ws.Range("E8").Value = '  +0.98%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.536'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.76%  '

$ws.Range("E10").Value = '  -0.97%  '

$ws.Range("E11").Value = '  -0.42%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.509'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +0.57%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.0000276'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.86%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '38.81'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +2.12%  '

$ws.Range("D15").Value = '3.741.41'
$ws.Range("E15").Value = '  +1.15%  '

$ws.Range("D16").Value = '66.191.10'
$ws.Range("E16").Value = '  +1.62%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '7.45'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +3.90%  '

$ws.Range("D18").Value = '3.221.19'

$ws.Range("E19").Value = '  +0.68%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '511.20'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '15.63'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +5.11%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.736'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.32%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '15.20'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.75%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '7.98'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +2.05%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '85.35'
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("E27").Value = '  +2.87%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '9.25'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +2.37%  '

$ws.Range("E29").Value = '  +2.63%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '6.90'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +9.53%  '

$ws.Range("E31").Value = '  +2.79%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '28.21'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +0.65%  '

$ws.Range("E33").Value = '  +1.21%  '

$ws.Range("E34").Value = '  +0.15%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '6.63'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +0.31%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '55.37'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.50%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.0915'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +1.75%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '483.88'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +1.38%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.0421'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.36%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.99'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -4.99%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '8.85'
$cell.Style = "Normal"

$ws.Range("E42").Value = '  +3.15%  '

$ws.Range("E43").Value = '  +0.49%  '

$ws.Range("E44").Value = '  +4.33%  '

$ws.Range("D45").Value = '2.959.33'
$ws.Range("E45").Value = '  -3.33%  '

$ws.Range("E46").Value = '  +5.39%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '28.94'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.69%  '

$ws.Range("E48").Value = '  +0.07%  '

$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("E50").Value = '  +2.90%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '34.04'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +5.85%  '
